$d = $word.ActiveDocument

# 1. Update the activation date.
$d.Content.Find.Execute("Ativação: 01/01/2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2022", 2)

# Helper: insert a new italic paragraph right after the paragraph that
# currently contains $anchorText, with text $newText.
function Add-ItalicParagraphAfter($anchorText, $newText) {
    $findRng = $d.Content
    $null = $findRng.Find.Execute($anchorText, $true, $false, $false, $false, $false,
                                   $true, 1, $false, "", 0)
    $para = $findRng.Paragraphs(1)
    $insertionPoint = $para.Range.End
    $null = $para.Range.InsertParagraphAfter()

    $probe = $d.Range($insertionPoint, $insertionPoint)
    $newPara = $probe.Paragraphs(1)
    $newPara.Range.Text = $newText

    $probe2 = $d.Range($insertionPoint, $insertionPoint)
    $newPara2 = $probe2.Paragraphs(1)
    $textOnly = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
    $textOnly.Font.Italic = 1
}

# 2. Add the English translation paragraphs. Done bottom-up so earlier
#    Find operations are not perturbed by paragraphs inserted further
#    down in the document.

# "Programa" section body -> add English translation paragraph after it.
Add-ItalicParagraphAfter `
    "1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente." `
    "- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment."

# "Programa resumido" section body -> add English translation paragraph after it.
Add-ItalicParagraphAfter `
    "1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão" `
    "The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes"

# "Objetivos" section body -> add English translation paragraph after it.
Add-ItalicParagraphAfter `
    "Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão" `
    "To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues."

# 3. Replace the whole Bibliografia paragraph content with the new reference
#    list. The old paragraph's first run carries xml:space="preserve", so
#    the text is cleared first and fresh text is inserted into the (now
#    empty) paragraph, which yields a plain, un-decorated run.
$bibRng = $d.Content
$null = $bibRng.Find.Execute("Gestão de Negócios: Visões e dimensões empresariais da oOrganização.", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
$bibPara = $bibRng.Paragraphs(1)
$bibStart = $bibPara.Range.Start
$bibEnd = $bibPara.Range.End - 1
$d.Range($bibStart, $bibEnd).Delete()
$bibProbe = $d.Range($bibStart, $bibStart)
$bibPara2 = $bibProbe.Paragraphs(1)
$bibPara2.Range.InsertAfter("LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.")
